$d = $word.ActiveDocument

$d.Content.Find.Execute("0.122", $true, $false, $false, $false, $false,
                         $true, 1, $false, "0.118", 2)

$d.Content.Find.Execute("[0.108, 0.139]", $true, $false, $false, $false, $false,
                         $true, 1, $false, "[0.111, 0.124]", 2)

$d.Content.Find.Execute("0.795", $true, $false, $false, $false, $false,
                         $true, 1, $false, "0.830", 2)

$d.Content.Find.Execute("0.666", $true, $false, $false, $false, $false,
                         $true, 1, $false, "0.745", 2)

$d.Content.Find.Execute("0.052", $true, $false, $false, $false, $false,
                         $true, 1, $false, "0.043", 2)

$d.Content.Find.Execute("561.52", $true, $false, $false, $false, $false,
                         $true, 1, $false, "1315.74", 2)
